# Update "想去人数" (want-to-go count) figures for a handful of events
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2074
$ws1.Range("F4").Value = 858
$ws1.Range("F5").Value = 1208
$ws1.Range("F6").Value = 357

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2074
$ws4.Range("F6").Value = 858
$ws4.Range("F7").Value = 1208
$ws4.Range("F8").Value = 357
